# "Fruta / hortaliza, semanal"
# A new weekly price record (2022-12-07, serial 44902) is inserted at the
# top of the Macroferia Regional de Talca / Arandano (blue) data block,
# pushing all existing records (previously rows 44-90) down by one row
# (new rows 45-91). Sheet used range grows from A1:T90 to A1:T91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the first data row of the block (row 44),
# shifting rows 44:90 down to 45:91 and carrying their formatting along.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record.
$ws.Cells.Item(44, 1).Value  = 5
$ws.Cells.Item(44, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(44, 3).Value  = "Maule"
$ws.Cells.Item(44, 4).Value  = 44902
$ws.Cells.Item(44, 5).Value  = 7
$ws.Cells.Item(44, 6).Value  = "Fruta"
$ws.Cells.Item(44, 7).Value  = 100101
$ws.Cells.Item(44, 8).Value  = "Berries"
$ws.Cells.Item(44, 9).Value  = 100101001
$ws.Cells.Item(44, 10).Value = "Arándano (blue)"
$ws.Cells.Item(44, 11).Value = "Sin especificar"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 250
$ws.Cells.Item(44, 14).Value = 3000
$ws.Cells.Item(44, 15).Value = 3200
$ws.Cells.Item(44, 16).Value = 3080
$ws.Cells.Item(44, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(44, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(44, 19).Value = 1540
$ws.Cells.Item(44, 20).Value = 2
